$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 712.3333
$ws.Range("I12").Value = 712.3333
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 712.3333
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -542.3333
$ws.Range("N12").ClearContents()

$ws.Range("H43").Value = 1840.875
$ws.Range("I43").Value = 1852.6666
$ws.Range("J43").Value = 1833.8
$ws.Range("K43").Value = 1852.6666
$ws.Range("L43").Value = 1833.8
$ws.Range("M43").Value = -1783.6666
$ws.Range("N43").Value = -1971.8

$ws.Range("H58").Value = 1445.4286
$ws.Range("J58").Value = 700
$ws.Range("L58").Value = 2100
$ws.Range("N58").Value = -2400

$ws.Range("H68").Value = 25295
$ws.Range("J68").Value = 25295
$ws.Range("L68").Value = 25295
$ws.Range("N68").Value = -26793

$ws.Range("H71").Value = 25295
$ws.Range("J71").Value = 25295
$ws.Range("L71").Value = 75885
$ws.Range("N71").Value = -83373

$ws.Range("H103").Value = 1605.8462
$ws.Range("I103").Value = 2142.4443
$ws.Range("J103").Value = 398.5
$ws.Range("K103").Value = 6427.3329
$ws.Range("L103").Value = 1195.5
$ws.Range("M103").Value = -5841.3329
$ws.Range("N103").Value = -2367.5

$ws.Range("H135").Value = 723
$ws.Range("I135").Value = 645.2143
$ws.Range("K135").Value = 5806.928699999999
$ws.Range("M135").Value = -3271.928699999999

$ws.Range("H138").Value = 2670.5334
$ws.Range("I138").Value = 2113.7827
$ws.Range("J138").Value = 4499.857
$ws.Range("K138").Value = 6341.348100000001
$ws.Range("L138").Value = 13499.571
$ws.Range("M138").Value = -1201.348100000001
$ws.Range("N138").Value = -23779.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3064.8948
$ws.Range("I32").Value = 3064.8948
$ws.Range("K32").Value = 3064.8948
$ws.Range("M32").Value = -2777.8948

$ws.Range("H64").Value = 100000
$ws.Range("J64").Value = 100000
$ws.Range("L64").Value = 100000
$ws.Range("N64").Value = -100496

$ws.Range("H67").Value = 100000
$ws.Range("J67").Value = 100000
$ws.Range("L67").Value = 100000
$ws.Range("N67").Value = -101716

$ws.Range("H74").Value = 114695.08
$ws.Range("I74").Value = 135596.42
$ws.Range("K74").Value = 135596.42
$ws.Range("M74").Value = -134722.42

$ws.Range("H77").Value = 114695.08
$ws.Range("I77").Value = 135596.42
$ws.Range("K77").Value = 677982.1000000001
$ws.Range("M77").Value = -673614.1000000001

$ws.Range("H97").Value = 577.129
$ws.Range("I97").Value = 688.125
$ws.Range("K97").Value = 688.125
$ws.Range("M97").Value = -192.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 66667136
$ws.Range("J80").Value = 520.7273
$ws.Range("L80").Value = 520.7273
$ws.Range("N80").Value = -2516.7273

$ws.Range("H83").Value = 66667136
$ws.Range("J83").Value = 520.7273
$ws.Range("L83").Value = 2603.6365
$ws.Range("N83").Value = -12587.6365

$ws.Range("H134").Value = 3178.074
$ws.Range("I134").Value = 2652.9333
$ws.Range("J134").Value = 3834.5
$ws.Range("K134").Value = 7958.7999
$ws.Range("L134").Value = 11503.5
$ws.Range("M134").Value = -5423.7999
$ws.Range("N134").Value = -16573.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5886.7856
$ws.Range("I31").Value = 4767.5713
$ws.Range("J31").Value = 7006
$ws.Range("K31").Value = 4767.5713
$ws.Range("L31").Value = 7006
$ws.Range("M31").Value = -4472.5713
$ws.Range("N31").Value = -7596

$ws.Range("H34").Value = 5886.7856
$ws.Range("I34").Value = 4767.5713
$ws.Range("J34").Value = 7006
$ws.Range("K34").Value = 4767.5713
$ws.Range("L34").Value = 7006
$ws.Range("M34").Value = -4565.5713
$ws.Range("N34").Value = -7410

$ws.Range("H86").Value = 7113.722
$ws.Range("I86").Value = 7445.5835
$ws.Range("J86").Value = 6450
$ws.Range("K86").Value = 7445.5835
$ws.Range("L86").Value = 6450
$ws.Range("M86").Value = -6322.5835
$ws.Range("N86").Value = -8696

$ws.Range("H89").Value = 7113.722
$ws.Range("I89").Value = 7445.5835
$ws.Range("J89").Value = 6450
$ws.Range("K89").Value = 37227.9175
$ws.Range("L89").Value = 32250
$ws.Range("M89").Value = -31611.9175
$ws.Range("N89").Value = -43482

$ws.Range("H122").Value = 1915
$ws.Range("I122").Value = 1975.6471
$ws.Range("J122").Value = 1399.5
$ws.Range("K122").Value = 5926.9413
$ws.Range("L122").Value = 4198.5
$ws.Range("M122").Value = -3476.9413
$ws.Range("N122").Value = -9098.5

$ws.Range("H132").Value = 11909857
$ws.Range("J132").Value = 33339984
$ws.Range("L132").Value = 100019952
$ws.Range("N132").Value = -100025012

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1051.1482
$ws.Range("I2").Value = 258.73334
$ws.Range("K2").Value = 1552.40004
$ws.Range("M2").Value = -1439.40004

$ws.Range("H17").Value = 3108.3333
$ws.Range("J17").Value = 3108.3333
$ws.Range("L17").Value = 9324.999899999999
$ws.Range("N17").Value = -9662.999899999999

$ws.Range("H34").Value = 593.75
$ws.Range("J34").Value = 1750
$ws.Range("L34").Value = 5250
$ws.Range("N34").Value = -5418

$ws.Range("H38").Value = 385.46155
$ws.Range("I38").Value = 51.5
$ws.Range("J38").Value = 919.8
$ws.Range("K38").Value = 154.5
$ws.Range("L38").Value = 2759.4
$ws.Range("M38").Value = 192.5
$ws.Range("N38").Value = -3453.4

$ws.Range("H39").Value = 2506
$ws.Range("J39").Value = 3195.5
$ws.Range("L39").Value = 9586.5
$ws.Range("N39").Value = -10174.5

$ws.Range("H47").Value = 2319.9
$ws.Range("J47").Value = 2928.5715
$ws.Range("L47").Value = 8785.7145
$ws.Range("N47").Value = -9647.7145

$ws.Range("H55").Value = 3857.4666
$ws.Range("J55").Value = 4130.1665
$ws.Range("L55").Value = 12390.4995
$ws.Range("N55").Value = -12744.4995

$ws.Range("H113").Value = 1262
$ws.Range("J113").Value = 1461.1333
$ws.Range("L113").Value = 4383.3999
$ws.Range("N113").Value = -8723.3999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 9958.416999999999
$ws.Range("I102").Value = 1268.5
$ws.Range("K102").Value = 1268.5
$ws.Range("M102").Value = 353.5

$ws.Range("H105").Value = 79830
$ws.Range("J105").Value = 79830
$ws.Range("L105").Value = 79830
$ws.Range("N105").Value = -86818

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3999
$ws.Range("I61").Value = 3999
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3999
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3797
$ws.Range("N61").ClearContents()

$ws.Range("H96").Value = 49999
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 49999
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 49999
$ws.Range("N96").Value = -55491
$ws.Range("M96").ClearContents()

$ws.Range("H113").Value = 3999
$ws.Range("I113").Value = 3999
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3999
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1829
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 12427.286
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H132").Value = 5333
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6706.8335
$ws.Range("I62").Value = 6586.25
$ws.Range("K62").Value = 6586.25
$ws.Range("M62").Value = -5962.25

$ws.Range("H65").Value = 6706.8335
$ws.Range("I65").Value = 6586.25
$ws.Range("K65").Value = 32931.25
$ws.Range("M65").Value = -29811.25

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

$ws.Range("H107").Value = 577.13336
$ws.Range("I107").Value = 320.8
$ws.Range("J107").Value = 1089.8
$ws.Range("K107").Value = 962.4000000000001
$ws.Range("L107").Value = 3269.4
$ws.Range("M107").Value = 957.5999999999999
$ws.Range("N107").Value = -7109.4

$ws.Range("H113").Value = 1370.5
$ws.Range("I113").Value = 1370.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4111.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1941.5
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 14479.75
$ws.Range("I132").Value = 22500
$ws.Range("J132").Value = 6459.5
$ws.Range("K132").Value = 67500
$ws.Range("L132").Value = 19378.5
$ws.Range("M132").Value = -64970
$ws.Range("N132").Value = -24438.5
